$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 11 (duplicate P04004 entry) - rows below shift up
$ws.Rows.Item(11).Delete()

# Update the selected cell to match the target state
$ws.Range("B17").Select()
